$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# --- Update the refreshed query timestamps on the "data" sheet (column F) ---
$data.Range("F2").Value = "2021-10-05 14:22:45.801090"
$data.Range("F3").Value = "2021-10-05 14:22:45.801098"
$data.Range("F4").Value = "2021-10-05 14:22:45.801101"
$data.Range("F5").Value = "2021-10-05 14:22:45.801104"
$data.Range("F6").Value = "2021-10-05 14:22:45.801107"
$data.Range("F7").Value = "2021-10-05 14:22:45.801111"
$data.Range("F8").Value = "2021-10-05 14:22:45.801114"
$data.Range("F9").Value = "2021-10-05 14:22:45.801116"
$data.Range("F10").Value = "2021-10-05 14:22:45.801119"

# --- Add a new "metadata" sheet right after "data" ---
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$meta.Name = "metadata"

# Header row (column labels)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Silver Russell syndrome"
$meta.Range("C2").Value = 199

# data_version must be the literal text "1.11" (not a number). Build it in a
# scratch cell formatted as text, then copy/paste just the value into D2 so
# D2 ends up as a plain text cell; the scratch column is then removed again.
$meta.Range("Z1").NumberFormat = "@"
$meta.Range("Z1").Value = [string]"1.11"
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Columns.Item(26).Delete()

$meta.Range("E2").Value = "2021-01-29T10:47:29.173646Z"
$meta.Range("F2").Value = "2021-10-05 14:22:45.797382"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/199/?format=json"

# Match the bold/bordered/centered header style used on the "data" sheet's
# header row, copying format only so no new style entries are introduced.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$meta.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = $false
$data.Select()
